{"js": "// Replace the date line and all 100 answer-table equations in document\n// order. Each paragraph in the document body (the standalone date\n// paragraph followed by one paragraph per table cell, read row-by-row,\n// left-to-right) gets its text swapped for the corresponding updated\n// value while leaving all run/paragraph formatting untouched (we use\n// insertText(..., \"Replace\") on the existing paragraph range instead of\n// rebuilding runs).\nconst newValues = [\n  \"2023-08-22 Tuesday\",\n  \"76-34=42\", \"85+9=94\", \"70-59=11\", \"23+57=80\", \"28+41=69\",\n  \"59+22=81\", \"74-58=16\", \"35-10=25\", \"13-9=4\", \"62-6=56\",\n  \"43+26=69\", \"79-30=49\", \"27+2=29\", \"34+50=84\", \"12+64=76\",\n  \"98-57=41\", \"93-64=29\", \"69-56=13\", \"3+46=49\", \"35+47=82\",\n  \"34+24=58\", \"49+23=72\", \"18-15=3\", \"87-1=86\", \"58-35=23\",\n  \"64+21=85\", \"46+4=50\", \"30+4=34\", \"77+21=98\", \"28+58=86\",\n  \"52-15=37\", \"52+26=78\", \"59+11=70\", \"71-15=56\", \"11-2=9\",\n  \"2+62=64\", \"53+11=64\", \"54-38=16\", \"26+24=50\", \"78+7=85\",\n  \"44+22=66\", \"11+23=34\", \"65-54=11\", \"74+13=87\", \"15-1=14\",\n  \"85-1=84\", \"40-21=19\", \"31+3=34\", \"21-16=5\", \"2+29=31\",\n  \"75-65=10\", \"19-16=3\", \"10+64=74\", \"46+45=91\", \"3+39=42\",\n  \"59+38=97\", \"73-62=11\", \"73-27=46\", \"99-77=22\", \"98-10=88\",\n  \"23+54=77\", \"91-25=66\", \"47-2=45\", \"37+37=74\", \"72+4=76\",\n  \"16+40=56\", \"89-88=1\", \"45+3=48\", \"94-51=43\", \"51-27=24\",\n  \"24-12=12\", \"92-49=43\", \"39-11=28\", \"82+14=96\", \"30+59=89\",\n  \"22+66=88\", \"24+70=94\", \"34-25=9\", \"91+5=96\", \"92-29=63\",\n  \"43-34=9\", \"83+12=95\", \"94-80=14\", \"54-11=43\", \"21+53=74\",\n  \"40+59=99\", \"18+0=18\", \"69-69=0\", \"99-16=83\", \"58-22=36\",\n  \"38+56=94\", \"58+30=88\", \"28+20=48\", \"34+7=41\", \"63+28=91\",\n  \"43-19=24\", \"13+57=70\", \"59+9=68\", \"73+15=88\", \"64-28=36\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: \" +\n      paragraphs.items.length +\n      \" (expected \" +\n      newValues.length +\n      \")\"\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const oldText = para.text;\n  const newText = newValues[i];\n  if (oldText !== newText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every equation in the 20x5 answers table to\n# the new values, in document order (row-major: row 1 col 1..5, row 2\n# col 1..5, ...). Setting Range.Text on the cell / paragraph leaves the\n# existing run/paragraph formatting (fonts, size, alignment) untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph\n$d.Paragraphs.Item(1).Range.Text = \"2023-08-22 Tuesday\"\n\n# 2) Answers table (20 rows x 5 columns), new values in row-major order\n$newValues = @(\n    \"76-34=42\", \"85+9=94\", \"70-59=11\", \"23+57=80\", \"28+41=69\",\n    \"59+22=81\", \"74-58=16\", \"35-10=25\", \"13-9=4\", \"62-6=56\",\n    \"43+26=69\", \"79-30=49\", \"27+2=29\", \"34+50=84\", \"12+64=76\",\n    \"98-57=41\", \"93-64=29\", \"69-56=13\", \"3+46=49\", \"35+47=82\",\n    \"34+24=58\", \"49+23=72\", \"18-15=3\", \"87-1=86\", \"58-35=23\",\n    \"64+21=85\", \"46+4=50\", \"30+4=34\", \"77+21=98\", \"28+58=86\",\n    \"52-15=37\", \"52+26=78\", \"59+11=70\", \"71-15=56\", \"11-2=9\",\n    \"2+62=64\", \"53+11=64\", \"54-38=16\", \"26+24=50\", \"78+7=85\",\n    \"44+22=66\", \"11+23=34\", \"65-54=11\", \"74+13=87\", \"15-1=14\",\n    \"85-1=84\", \"40-21=19\", \"31+3=34\", \"21-16=5\", \"2+29=31\",\n    \"75-65=10\", \"19-16=3\", \"10+64=74\", \"46+45=91\", \"3+39=42\",\n    \"59+38=97\", \"73-62=11\", \"73-27=46\", \"99-77=22\", \"98-10=88\",\n    \"23+54=77\", \"91-25=66\", \"47-2=45\", \"37+37=74\", \"72+4=76\",\n    \"16+40=56\", \"89-88=1\", \"45+3=48\", \"94-51=43\", \"51-27=24\",\n    \"24-12=12\", \"92-49=43\", \"39-11=28\", \"82+14=96\", \"30+59=89\",\n    \"22+66=88\", \"24+70=94\", \"34-25=9\", \"91+5=96\", \"92-29=63\",\n    \"43-34=9\", \"83+12=95\", \"94-80=14\", \"54-11=43\", \"21+53=74\",\n    \"40+59=99\", \"18+0=18\", \"69-69=0\", \"99-16=83\", \"58-22=36\",\n    \"38+56=94\", \"58+30=88\", \"28+20=48\", \"34+7=41\", \"63+28=91\",\n    \"43-19=24\", \"13+57=70\", \"59+9=68\", \"73+15=88\", \"64-28=36\"\n)\n\n$tbl = $d.Tables.Item(1)\n$numRows = $tbl.Rows.Count\n$numCols = $tbl.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $numRows; $r++) {\n    for ($c = 1; $c -le $numCols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i = $i + 1\n    }\n}\n"}
